$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("testdata_Mean")
$ws2 = $wb.Worksheets.Item("testdata_Mean_results")

# The sheet currently ends with two blank (but style-formatted) rows, 28
# and 29. We want to fill those two rows with a "group total" style
# summary (Area1/NA and Area2/NA) and push a brand new blank row down to
# row 30 so the sheet still ends with one blank formatted row.
#
# Duplicate the still-empty row 29 (which already carries the correct
# cell formatting) and insert it as the new row 30, shifting nothing else
# - this keeps the formatting on the new trailing blank row instead of it
# coming out as a completely bare/default cell.
$ws1.Rows.Item(29).Copy()
$ws1.Rows.Item(30).Insert(-4121)   # -4121 = xlShiftDown

# Fill in the now-vacated rows 28 and 29 with the new summary data.
$ws1.Range("A28").Value = "Area1"
$ws1.Range("B28").Value = "NA"
$ws1.Range("A29").Value = "Area2"
$ws1.Range("B29").Value = "NA"

# testdata_Mean becomes the active sheet/tab, with the newly added B28
# cell selected (testdata_Mean_results loses its "tabSelected" flag as a
# result, matching the edit).
$ws1.Activate()
$ws1.Range("B28").Select()
